$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 779.7692
$ws.Range("I43").Value = 601.2
$ws.Range("J43").Value = 891.375
$ws.Range("K43").Value = 601.2
$ws.Range("L43").Value = 891.375
$ws.Range("M43").Value = -532.2
$ws.Range("N43").Value = -1029.375
$ws.Range("H98").Value = 8076405
$ws.Range("I98").Value = 11010552
$ws.Range("K98").Value = 11010552
$ws.Range("M98").Value = -11009054
$ws.Range("H122").Value = 8076405
$ws.Range("I122").Value = 11010552
$ws.Range("K122").Value = 33031656
$ws.Range("M122").Value = -33029206
$ws.Range("H129").Value = 897.47
$ws.Range("I129").Value = 313.63635
$ws.Range("J129").Value = 969.6292
$ws.Range("K129").Value = 940.90905
$ws.Range("L129").Value = 2908.8876
$ws.Range("M129").Value = 4059.09095
$ws.Range("N129").Value = -12908.8876
$ws.Range("H132").Value = 4599.552
$ws.Range("I132").Value = 4687.923
$ws.Range("J132").Value = 3833.6667
$ws.Range("K132").Value = 14063.769
$ws.Range("L132").Value = 11501.0001
$ws.Range("M132").Value = -11533.769
$ws.Range("N132").Value = -16561.0001
$ws.Range("H137").Value = 1569.9656
$ws.Range("I137").Value = 1611.9584
$ws.Range("J137").Value = 1368.4
$ws.Range("K137").Value = 4835.8752
$ws.Range("L137").Value = 4105.200000000001
$ws.Range("M137").Value = -2285.8752
$ws.Range("N137").Value = -9205.200000000001
$ws.Range("H138").Value = 3894.3914
$ws.Range("I138").Value = 1927.3334
$ws.Range("J138").Value = 5546.72
$ws.Range("K138").Value = 5782.0002
$ws.Range("L138").Value = 16640.16
$ws.Range("M138").Value = -642.0002000000004
$ws.Range("N138").Value = -26920.16
$ws.Range("H140").Value = 74593.336
$ws.Range("J140").Value = 74593.336
$ws.Range("L140").Value = 74593.336
$ws.Range("N140").Value = -84953.336

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2641.3
$ws.Range("I45").Value = 2321.25
$ws.Range("K45").Value = 2321.25
$ws.Range("M45").Value = -1944.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 473
$ws.Range("I64").Value = 459.33334
$ws.Range("J64").Value = 479.83334
$ws.Range("K64").Value = 459.33334
$ws.Range("L64").Value = 479.83334
$ws.Range("M64").Value = -234.33334
$ws.Range("N64").Value = -929.83334
$ws.Range("H67").Value = 473
$ws.Range("I67").Value = 459.33334
$ws.Range("J67").Value = 479.83334
$ws.Range("K67").Value = 459.33334
$ws.Range("L67").Value = 479.83334
$ws.Range("M67").Value = 320.66666
$ws.Range("N67").Value = -2039.83334
$ws.Range("H82").Value = 13900.363
$ws.Range("I82").Value = 2180.25
$ws.Range("J82").Value = 20597.572
$ws.Range("K82").Value = 2180.25
$ws.Range("L82").Value = 20597.572
$ws.Range("M82").Value = -1797.25
$ws.Range("N82").Value = -21363.572
$ws.Range("H85").Value = 13900.363
$ws.Range("I85").Value = 2180.25
$ws.Range("J85").Value = 20597.572
$ws.Range("K85").Value = 2180.25
$ws.Range("L85").Value = 20597.572
$ws.Range("M85").Value = -854.25
$ws.Range("N85").Value = -23249.572
$ws.Range("H107").Value = 1462
$ws.Range("I107").Value = 1194.125
$ws.Range("J107").Value = 1700.1111
$ws.Range("K107").Value = 1194.125
$ws.Range("L107").Value = 1700.1111
$ws.Range("M107").Value = 725.875
$ws.Range("N107").Value = -5540.1111
$ws.Range("H134").Value = 7450.5
$ws.Range("I134").Value = 2513.0908
$ws.Range("J134").Value = 21028.375
$ws.Range("K134").Value = 7539.2724
$ws.Range("L134").Value = 63085.125
$ws.Range("M134").Value = -5004.2724
$ws.Range("N134").Value = -68155.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1025.4286
$ws.Range("I16").Value = 935.85
$ws.Range("J16").Value = 1249.375
$ws.Range("K16").Value = 935.85
$ws.Range("L16").Value = 1249.375
$ws.Range("M16").Value = -648.85
$ws.Range("N16").Value = -1823.375
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").ClearContents()
$ws.Range("H31").Value = 5210159
$ws.Range("I31").Value = 6945868.5
$ws.Range("J31").Value = 3030.3333
$ws.Range("K31").Value = 6945868.5
$ws.Range("L31").Value = 3030.3333
$ws.Range("M31").Value = -6945573.5
$ws.Range("N31").Value = -3620.3333
$ws.Range("H34").Value = 5210159
$ws.Range("I34").Value = 6945868.5
$ws.Range("J34").Value = 3030.3333
$ws.Range("K34").Value = 6945868.5
$ws.Range("L34").Value = 3030.3333
$ws.Range("M34").Value = -6945666.5
$ws.Range("N34").Value = -3434.3333
$ws.Range("H113").Value = 1025.4286
$ws.Range("I113").Value = 935.85
$ws.Range("J113").Value = 1249.375
$ws.Range("K113").Value = 935.85
$ws.Range("L113").Value = 1249.375
$ws.Range("M113").Value = 1234.15
$ws.Range("N113").Value = -5589.375
$ws.Range("H132").Value = 38121.645
$ws.Range("I132").Value = 1857.238
$ws.Range("J132").Value = 146914.86
$ws.Range("K132").Value = 5571.714
$ws.Range("L132").Value = 440744.58
$ws.Range("M132").Value = -3041.714
$ws.Range("N132").Value = -445804.58
$ws.Range("H134").Value = 2746.28
$ws.Range("I134").Value = 1670.1364
$ws.Range("J134").Value = 10638
$ws.Range("K134").Value = 5010.4092
$ws.Range("L134").Value = 31914
$ws.Range("M134").Value = -2475.4092
$ws.Range("N134").Value = -36984

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 524.36365
$ws.Range("J92").Value = 549.3333
$ws.Range("L92").Value = 1647.9999
$ws.Range("N92").Value = -4143.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1448.4
$ws.Range("I122").Value = 1340.2
$ws.Range("J122").Value = 1881.2
$ws.Range("K122").Value = 4020.6
$ws.Range("L122").Value = 5643.6
$ws.Range("M122").Value = -1570.6
$ws.Range("N122").Value = -10543.6
$ws.Range("H132").Value = 334774.1
$ws.Range("I132").Value = 41109.32
$ws.Range("J132").Value = 1252476.5
$ws.Range("K132").Value = 123327.96
$ws.Range("L132").Value = 3757429.5
$ws.Range("M132").Value = -120797.96
$ws.Range("N132").Value = -3762489.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 6985.9414
$ws.Range("I46").Value = 880.1
$ws.Range("J46").Value = 15708.571
$ws.Range("K46").Value = 880.1
$ws.Range("L46").Value = 15708.571
$ws.Range("M46").Value = -692.1
$ws.Range("N46").Value = -16084.571
